$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.487.77"
$ws.Range("E2").Value = "  -3.10%  "

$ws.Range("D3").Value = "3.725.68"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'591.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.54%  "

$ws.Range("D6").Value = "'171.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.25%  "

$ws.Range("D7").Value = "3.725.53"
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("E10").Value = "  -5.44%  "

$ws.Range("D11").Value = "'6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.00%  "

$ws.Range("D12").Value = "'0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.88%  "

$ws.Range("D13").Value = "'37.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.79%  "

$ws.Range("D14").Value = "'0.0000240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.19%  "

$ws.Range("D15").Value = "4.336.80"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").Value = "3.715.99"
$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "67.396.06"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("E18").Value = "  -5.30%  "

$ws.Range("D19").Value = "'7.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.52%  "

$ws.Range("D20").Value = "'16.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").Value = "'485.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").Value = "'8.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.47%  "

$ws.Range("D23").Value = "'0.710"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.53%  "

$ws.Range("D24").Value = "'83.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.14%  "

$ws.Range("D25").Value = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.46%  "

$ws.Range("D26").Value = "'0.0000139"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.63%  "

$ws.Range("E27").Value = "  -6.22%  "

$ws.Range("D28").Value = "'10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.58%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'2.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'32.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.76%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.53%  "

$ws.Range("D33").Value = "'7.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.07%  "

$ws.Range("E34").Value = "  -5.24%  "

$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").Value = "'0.993"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.43%  "

$ws.Range("E37").Value = "  -2.91%  "

$ws.Range("D38").Value = "'5.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.31%  "

$ws.Range("D39").Value = "'0.321"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.18%  "

$ws.Range("D40").Value = "'447.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("D41").Value = "'48.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.91%  "

$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("D43").Value = "'2.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.86%  "

$ws.Range("E44").Value = "  -4.06%  "

$ws.Range("D45").Value = "'41.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.20%  "

$ws.Range("D46").Value = "'140.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("D47").Value = "2.780.51"
$ws.Range("E47").Value = "  -5.85%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "'0.0345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.94%  "

$ws.Range("D50").Value = "'25.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.47%  "

$ws.Range("D51").Value = "'22.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.43%  "
